$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6297.5
$ws.Range("I62").Value = 3652
$ws.Range("J62").Value = 8187.143
$ws.Range("K62").Value = 3652
$ws.Range("L62").Value = 8187.143
$ws.Range("M62").Value = -3028
$ws.Range("N62").Value = -9435.143

$ws.Range("H65").Value = 6297.5
$ws.Range("I65").Value = 3652
$ws.Range("J65").Value = 8187.143
$ws.Range("K65").Value = 18260
$ws.Range("L65").Value = 40935.715
$ws.Range("M65").Value = -15140
$ws.Range("N65").Value = -47175.715

$ws.Range("H98").Value = 699838.7
$ws.Range("I98").Value = 799208.5
$ws.Range("J98").Value = 4250
$ws.Range("K98").Value = 799208.5
$ws.Range("L98").Value = 4250
$ws.Range("M98").Value = -797710.5
$ws.Range("N98").Value = -7246

$ws.Range("H122").Value = 699838.7
$ws.Range("I122").Value = 799208.5
$ws.Range("J122").Value = 4250
$ws.Range("K122").Value = 2397625.5
$ws.Range("L122").Value = 12750
$ws.Range("M122").Value = -2395175.5
$ws.Range("N122").Value = -17650

$ws.Range("H138").Value = 6191240.5
$ws.Range("I138").Value = 2305228.2
$ws.Range("J138").Value = 6946853.5
$ws.Range("K138").Value = 6915684.600000001
$ws.Range("L138").Value = 20840560.5
$ws.Range("M138").Value = -6910544.600000001
$ws.Range("N138").Value = -20850840.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5558.077
$ws.Range("I32").Value = 3730.353
$ws.Range("K32").Value = 3730.353
$ws.Range("M32").Value = -3443.353

$ws.Range("H37").Value = 4300
$ws.Range("J37").Value = 4300
$ws.Range("L37").Value = 4300
$ws.Range("N37").Value = -4846

$ws.Range("H110").Value = 1334.5555
$ws.Range("I110").Value = 1388.875
$ws.Range("J110").Value = 900
$ws.Range("K110").Value = 1388.875
$ws.Range("L110").Value = 900
$ws.Range("M110").Value = 656.125
$ws.Range("N110").Value = -4990

$ws.Range("H122").Value = 1338.8572
$ws.Range("I122").Value = 1078.4
$ws.Range("J122").Value = 1990
$ws.Range("K122").Value = 3235.2
$ws.Range("L122").Value = 5970
$ws.Range("M122").Value = -785.2000000000003
$ws.Range("N122").Value = -10870

$ws.Range("H132").Value = 2398.7585
$ws.Range("I132").Value = 1972.9565
$ws.Range("J132").Value = 4031
$ws.Range("K132").Value = 5918.8695
$ws.Range("L132").Value = 12093
$ws.Range("M132").Value = -3388.8695
$ws.Range("N132").Value = -17153

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2601.5417
$ws.Range("I132").Value = 1913.3846
$ws.Range("J132").Value = 3414.818
$ws.Range("K132").Value = 5740.1538
$ws.Range("L132").Value = 10244.454
$ws.Range("M132").Value = -3210.1538
$ws.Range("N132").Value = -15304.454

$ws.Range("H134").Value = 3014.1667
$ws.Range("I134").Value = 1375.8
$ws.Range("J134").Value = 4184.4287
$ws.Range("K134").Value = 4127.4
$ws.Range("L134").Value = 12553.2861
$ws.Range("M134").Value = -1592.4
$ws.Range("N134").Value = -17623.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 866.6667
$ws.Range("J46").Value = 866.6667
$ws.Range("L46").Value = 2600.0001
$ws.Range("N46").Value = -2782.0001

$ws.Range("H51").Value = 919.875
$ws.Range("I51").Value = 702
$ws.Range("J51").Value = 992.5
$ws.Range("K51").Value = 2106
$ws.Range("L51").Value = 2977.5
$ws.Range("M51").Value = -1646
$ws.Range("N51").Value = -3897.5

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()

$ws.Range("H117").Value = 1539
$ws.Range("I117").Value = 379.4
$ws.Range("J117").Value = 2505.3333
$ws.Range("K117").Value = 1138.2
$ws.Range("L117").Value = 7515.999899999999
$ws.Range("M117").Value = 2303.8
$ws.Range("N117").Value = -14399.9999

$ws.Range("H129").Value = 5000
$ws.Range("I129").Value = 10000
$ws.Range("J129").Value = 2500
$ws.Range("K129").Value = 30000
$ws.Range("L129").Value = 7500
$ws.Range("M129").Value = -25000
$ws.Range("N129").Value = -17500

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 307
$ws.Range("J107").Value = 361
$ws.Range("L107").Value = 361
$ws.Range("N107").Value = -4201

$ws.Range("H122").Value = 3704406
$ws.Range("I122").Value = 3704406
$ws.Range("K122").Value = 11113218
$ws.Range("M122").Value = -11110768

$ws.Range("H126").Value = 2535.7144
$ws.Range("I126").Value = 2750
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 8250
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -5780
$ws.Range("N126").Value = -12440

$ws.Range("H132").Value = 2020.9231
$ws.Range("I132").Value = 1677.2667
$ws.Range("J132").Value = 3166.4443
$ws.Range("K132").Value = 5031.800099999999
$ws.Range("L132").Value = 9499.332900000001
$ws.Range("M132").Value = -2501.800099999999
$ws.Range("N132").Value = -14559.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H74").Value = 98599.14
$ws.Range("I74").Value = 275197
$ws.Range("J74").Value = 27960
$ws.Range("K74").Value = 275197
$ws.Range("L74").Value = 27960
$ws.Range("M74").Value = -274199
$ws.Range("N74").Value = -29956

$ws.Range("H77").Value = 98599.14
$ws.Range("I77").Value = 275197
$ws.Range("J77").Value = 27960
$ws.Range("K77").Value = 825591
$ws.Range("L77").Value = 83880
$ws.Range("M77").Value = -820599
$ws.Range("N77").Value = -93864

$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -32246

$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -101232

$ws.Range("H132").Value = 4499.625
$ws.Range("I132").Value = 5000
$ws.Range("J132").Value = 4428.143
$ws.Range("K132").Value = 15000
$ws.Range("L132").Value = 13284.429
$ws.Range("M132").Value = -12470
$ws.Range("N132").Value = -18344.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 51000
$ws.Range("I57").Value = 51000
$ws.Range("K57").Value = 51000
$ws.Range("M57").Value = -50246

$ws.Range("H132").Value = 13515560
$ws.Range("I132").Value = 18520072
$ws.Range("J132").Value = 3378.3
$ws.Range("K132").Value = 55560216
$ws.Range("L132").Value = 10134.9
$ws.Range("M132").Value = -55557686
$ws.Range("N132").Value = -15194.9

$ws.Range("H136").Value = 13931756
$ws.Range("I136").Value = 15921664
$ws.Range("J136").Value = 2398.3333
$ws.Range("K136").Value = 47764992
$ws.Range("L136").Value = 7194.999899999999
$ws.Range("M136").Value = -47762442
$ws.Range("N136").Value = -12294.9999
